$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 298.5
$ws.Range("I6").Value = 100
$ws.Range("J6").Value = 364.66666
$ws.Range("K6").Value = 300
$ws.Range("L6").Value = 1093.99998
$ws.Range("M6").Value = -188
$ws.Range("N6").Value = -1317.99998
$ws.Range("H106").Value = 1295.3334
$ws.Range("I106").Value = 1295.3334
$ws.Range("K106").Value = 1295.3334
$ws.Range("M106").Value = -664.3334
$ws.Range("H138").Value = 3598.146
$ws.Range("I138").Value = 2334.9048
$ws.Range("J138").Value = 3988.2646
$ws.Range("K138").Value = 7004.714399999999
$ws.Range("L138").Value = 11964.7938
$ws.Range("M138").Value = -1864.714399999999
$ws.Range("N138").Value = -22244.7938
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 65875.44
$ws.Range("I74").Value = 72517.72
$ws.Range("J74").Value = 18430.572
$ws.Range("K74").Value = 72517.72
$ws.Range("L74").Value = 18430.572
$ws.Range("M74").Value = -71643.72
$ws.Range("N74").Value = -20178.572
$ws.Range("H77").Value = 65875.44
$ws.Range("I77").Value = 72517.72
$ws.Range("J77").Value = 18430.572
$ws.Range("K77").Value = 362588.6
$ws.Range("L77").Value = 92152.86
$ws.Range("M77").Value = -358220.6
$ws.Range("N77").Value = -100888.86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 6499.25
$ws.Range("J15").Value = 6499.25
$ws.Range("L15").Value = 6499.25
$ws.Range("N15").Value = -6953.25
$ws.Range("H105").Value = 4959.0312
$ws.Range("I105").Value = 4590.8335
$ws.Range("K105").Value = 4590.8335
$ws.Range("M105").Value = -2843.8335
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 10093.23
$ws.Range("J4").Value = 10093.23
$ws.Range("L4").Value = 10093.23
$ws.Range("N4").Value = -10317.23
$ws.Range("H14").Value = 700
$ws.Range("I14").Value = 700
$ws.Range("K14").Value = 700
$ws.Range("M14").Value = -530
$ws.Range("H16").Value = 1918.7778
$ws.Range("I16").Value = 1924.1428
$ws.Range("J16").Value = 1900
$ws.Range("K16").Value = 1924.1428
$ws.Range("L16").Value = 1900
$ws.Range("M16").Value = -1637.1428
$ws.Range("N16").Value = -2474
$ws.Range("H113").Value = 1918.7778
$ws.Range("I113").Value = 1924.1428
$ws.Range("J113").Value = 1900
$ws.Range("K113").Value = 1924.1428
$ws.Range("L113").Value = 1900
$ws.Range("M113").Value = 245.8571999999999
$ws.Range("N113").Value = -6240
$ws.Range("H122").Value = 7629.2812
$ws.Range("I122").Value = 4343.077
$ws.Range("J122").Value = 9877.736999999999
$ws.Range("K122").Value = 13029.231
$ws.Range("L122").Value = 29633.211
$ws.Range("M122").Value = -10579.231
$ws.Range("N122").Value = -34533.211
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 934.2162
$ws.Range("I107").Value = 389.66666
$ws.Range("K107").Value = 1168.99998
$ws.Range("M107").Value = 751.0000199999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 14000
$ws.Range("J12").Value = 14000
$ws.Range("L12").Value = 14000
$ws.Range("N12").Value = -14280
$ws.Range("H102").Value = 2603.0322
$ws.Range("I102").Value = 2318
$ws.Range("K102").Value = 2318
$ws.Range("M102").Value = -696
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3669.2856
$ws.Range("I7").Value = 3407.2727
$ws.Range("J7").Value = 3957.5
$ws.Range("K7").Value = 3407.2727
$ws.Range("L7").Value = 3957.5
$ws.Range("M7").Value = -3295.2727
$ws.Range("N7").Value = -4181.5
$ws.Range("H16").Value = 1789
$ws.Range("I16").Value = 1683.5
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1683.5
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -1513.5
$ws.Range("N16").Value = -2340
$ws.Range("H17").Value = 933.3333
$ws.Range("I17").Value = 933.3333
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 933.3333
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -763.3333
$ws.Range("N17").ClearContents()
$ws.Range("H61").Value = 2014268.4
$ws.Range("I61").Value = 2756130.2
$ws.Range("J61").Value = 35970
$ws.Range("K61").Value = 2756130.2
$ws.Range("L61").Value = 35970
$ws.Range("M61").Value = -2755928.2
$ws.Range("N61").Value = -36374
$ws.Range("H113").Value = 2014268.4
$ws.Range("I113").Value = 2756130.2
$ws.Range("J113").Value = 35970
$ws.Range("K113").Value = 2756130.2
$ws.Range("L113").Value = 35970
$ws.Range("M113").Value = -2753960.2
$ws.Range("N113").Value = -40310
$ws.Range("H126").Value = 3669.2856
$ws.Range("I126").Value = 3407.2727
$ws.Range("J126").Value = 3957.5
$ws.Range("K126").Value = 10221.8181
$ws.Range("L126").Value = 11872.5
$ws.Range("M126").Value = -7751.8181
$ws.Range("N126").Value = -16812.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 2000
$ws.Range("J5").Value = 2000
$ws.Range("L5").Value = 2000
$ws.Range("N5").Value = -2224
$ws.Range("H14").Value = 5500
$ws.Range("I14").Value = 3250
$ws.Range("J14").Value = 10000
$ws.Range("K14").Value = 3250
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = -3082
$ws.Range("N14").Value = -10336
$ws.Range("H19").Value = 10000
$ws.Range("J19").Value = 10000
$ws.Range("L19").Value = 10000
$ws.Range("N19").Value = -10348
$ws.Range("H30").Value = 10450
$ws.Range("J30").Value = 9333.333000000001
$ws.Range("L30").Value = 9333.333000000001
$ws.Range("N30").Value = -9547.333000000001
$ws.Range("H33").Value = 19830.25
$ws.Range("J33").Value = 19830.25
$ws.Range("L33").Value = 19830.25
$ws.Range("N33").Value = -20330.25
$ws.Range("H36").Value = 19830.25
$ws.Range("J36").Value = 19830.25
$ws.Range("L36").Value = 19830.25
$ws.Range("N36").Value = -20330.25
$ws.Range("H122").Value = 3564.2942
$ws.Range("I122").Value = 2015.3529
$ws.Range("J122").Value = 5113.2354
$ws.Range("K122").Value = 6046.0587
$ws.Range("L122").Value = 15339.7062
$ws.Range("M122").Value = -3596.0587
$ws.Range("N122").Value = -20239.7062
$ws.Range("H124").Value = 56219.5
$ws.Range("J124").Value = 56219.5
$ws.Range("L124").Value = 56219.5
$ws.Range("N124").Value = -66039.5
$ws.Range("H126").Value = 1916.6666
$ws.Range("I126").Value = 1900
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 5700
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -3230
$ws.Range("N126").Value = -10940
$ws.Range("H129").Value = 49450
$ws.Range("J129").Value = 49450
$ws.Range("L129").Value = 49450
$ws.Range("N129").Value = -59450
$ws.Range("H137").Value = 58068
$ws.Range("J137").Value = 58068
$ws.Range("L137").Value = 58068
$ws.Range("N137").Value = -68268
